$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 of data, matching style of rows 3-5 (style index "1": thin borders all around)
$ws.Range("A6").Value = "W3"

$values = @(8.7200000000000006, 8.42, 9.36, 11.09, 9.75, 8.34, 9.86, 10.88, 11.27, 10.56, 12.78, 11.14)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i  # column B is 2
    $cell = $ws.Cells.Item(6, $col)
    $cell.Value = $values[$i]
}

# Copy style (borders) from row 5 (B5:M5) to row 6 (B6:M6) so it matches style "1"
$ws.Range("B5:M5").Copy()
$ws.Range("B6:M6").PasteSpecial(-4122) # xlPasteFormats

# Set the selection to E8 as per the diff
$ws.Range("E8").Select()
